# Auto-generated edit script: updates cryptocurrency price/volume data
# to reflect the latest GitHub Actions scrape, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number plus the new values for any of columns B (Coin),
# C (Link), D (Price) and E (Volume(1h)) that changed for that row.
# Columns not present for a row are left untouched.
$updates = @(
    @{Row=2; D='35.290.58'; E='  +0.30%  '},
    @{Row=3; D='1.879.70'; E='  -1.30%  '},
    @{Row=4; E='  -0.52%  '},
    @{Row=5; D='246.72'; E='  -2.85%  '},
    @{Row=6; E='  -2.34%  '},
    @{Row=7; E='  -0.58%  '},
    @{Row=8; D='43.59'; E='  +4.34%  '},
    @{Row=9; E='  +0.26%  '},
    @{Row=10; D='53.70'; E='  +1.73%  '},
    @{Row=11; D='0.0740'; E='  -2.69%  '},
    @{Row=12; E='  -0.45%  '},
    @{Row=13; D='13.58'; E='  +2.95%  '},
    @{Row=14; D='2.152.40'; E='  -1.34%  '},
    @{Row=15; D='0.765'; E='  +3.69%  '},
    @{Row=16; E='  -2.07%  '},
    @{Row=17; D='1.874.08'; E='  -1.42%  '},
    @{Row=18; D='35.335.13'; E='  +0.48%  '},
    @{Row=19; D='72.75'; E='  -1.53%  '},
    @{Row=20; D='0.0₃0822'; E='  -2.69%  '},
    @{Row=21; D='244.09'; E='  +0.37%  '},
    @{Row=22; D='12.84'; E='  -1.88%  '},
    @{Row=23; D='4.98'; E='  -1.67%  '},
    @{Row=24; E='  +7.37%  '},
    @{Row=25; E='  -0.52%  '},
    @{Row=26; D='2.20'; E='  -5.90%  '},
    @{Row=27; D='165.50'; E='  -1.52%  '},
    @{Row=28; D='8.57'; E='  -0.45%  '},
    @{Row=29; E='  -1.57%  '},
    @{Row=30; E='  -2.30%  '},
    @{Row=31; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.71'; E='  +4.83%  '},
    @{Row=32; B='WEMIXToken'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='2.04'; E='  +0.96%  '},
    @{Row=33; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.30'; E='  -1.08%  '},
    @{Row=34; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0593'; E='  -1.25%  '},
    @{Row=35; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.15'; E='  -2.41%  '},
    @{Row=36; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.00'; E='  -0.53%  '},
    @{Row=37; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.838'; E='  -2.00%  '},
    @{Row=38; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.94'; E='  -3.61%  '},
    @{Row=39; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.0718'; E='  +9.14%  '},
    @{Row=40; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='17.64'; E='  +2.84%  '},
    @{Row=41; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0218'; E='  +0.80%  '},
    @{Row=42; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='96.35'; E='  -2.41%  '},
    @{Row=43; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.08'; E='  -3.15%  '},
    @{Row=44; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.304.34'; E='  -0.27%  '},
    @{Row=45; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.37'; E='  -2.17%  '},
    @{Row=46; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.0803'; E='  +6.47%  '},
    @{Row=47; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.38'; E='  -1.64%  '},
    @{Row=48; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.73'; E='  -0.75%  '},
    @{Row=49; B='Gas'; C='https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'; D='11.93'; E='  -2.91%  '},
    @{Row=50; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='6.23'; E='  -5.77%  '},
    @{Row=51; B='MultiversX'; C='https://coinranking.com/coin/omwkOTglq+multiversx-egld'; D='42.06'; E='  -2.39%  '}
)

# The Price column (D) holds values like '35.290.58' or '1.00' that Excel
# would otherwise reinterpret as numbers/dates, so force it to Text format
# before writing so the values round-trip exactly as strings.
$ws.Range("D2:D51").NumberFormat = "@"

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = $u.D }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
}
